$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loot_All")

# Row 10 - Pistol (entered ID, Name, then Path)
$ws.Cells.Item(10, 1).Value = "PISTOL_1"
$ws.Cells.Item(10, 3).Value = "Pistol"
$ws.Cells.Item(10, 2).Value = "Assets/Scripts/Weapons/Pistol.prefab"
$ws.Cells.Item(10, 4).Value = "COMMON"

# Rows 11-12 - Machine Gun / Shot Gun (entered column by column: ID, Path, Name)
$ws.Cells.Item(11, 1).Value = "MACHINE_GUN_1"
$ws.Cells.Item(12, 1).Value = "SHOT_GUN_1"

$ws.Cells.Item(11, 2).Value = "Assets/Scripts/Weapons/MachineGun.prefab"
$ws.Cells.Item(12, 2).Value = "Assets/Scripts/Weapons/ShotGun.prefab"

$ws.Cells.Item(11, 3).Value = "Machine Gun"
$ws.Cells.Item(12, 3).Value = "Shot Gun"

$ws.Cells.Item(11, 4).Value = "COMMON"
$ws.Cells.Item(12, 4).Value = "COMMON"

$ws.Range("D17").Select()
